$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "EDEN.PA"
$ws.Range("B14").Value = "Edenred"
$ws.Range("C14").Value = "Industrials"
$ws.Range("D14").Value = "Transaction Processing Services"

$ws.Range("A15").Value = "ENGI.PA"
$ws.Range("B15").Value = "Engie"
$ws.Range("C15").Value = "Utilities"
$ws.Range("D15").Value = "Gas Utilities"

$ws.Range("A16").Value = "EL.PA"
$ws.Range("B16").Value = "EssilorLuxottica"
$ws.Range("C16").Value = "Healthcare"
$ws.Range("D16").Value = "Apparel, Accessories & Luxury Goods"

$ws.Range("A17").Value = "ERF.PA"
$ws.Range("B17").Value = "Eurofins Scientific"
$ws.Range("C17").Value = "Healthcare"
$ws.Range("D17").Value = "Biotechnologies"

$ws.Range("A18").Value = "RMS.PA"
$ws.Range("B18").Value = "Hermès"
$ws.Range("C18").Value = "Consumer Cyclical"
$ws.Range("D18").Value = "Apparel, Accessories & Luxury Goods"

$ws.Range("A19").Value = "KER.PA"
$ws.Range("B19").Value = "Kering"
$ws.Range("C19").Value = "Consumer Cyclical"
$ws.Range("D19").Value = "Apparel, Accessories & Luxury Goods"

$ws.Range("A20").Value = "OR.PA"
$ws.Range("B20").Value = "L'Oréal"
$ws.Range("C20").Value = "Consumer Defensive"
$ws.Range("D20").Value = "Personal Products"

$ws.Range("A21").Value = "LR.PA"
$ws.Range("B21").Value = "Legrand"
$ws.Range("C21").Value = "Industrials"
$ws.Range("D21").Value = "Electrical Components & Equipment"

$ws.Range("A22").Value = "MC.PA"
$ws.Range("B22").Value = "LVMH"
$ws.Range("C22").Value = "Consumer Cyclical"
$ws.Range("D22").Value = "Apparel, Accessories & Luxury Goods"

$ws.Range("A23").Value = "ML.PA"
$ws.Range("B23").Value = "Michelin"
$ws.Range("C23").Value = "Industrials"
$ws.Range("D23").Value = "Tires & Rubber"

$ws.Range("A24").Value = "ORA.PA"
$ws.Range("B24").Value = "Orange"
$ws.Range("C24").Value = "Communication Services"
$ws.Range("D24").Value = "Integrated Telecommunication Services"

$ws.Range("A25").Value = "RI.PA"
$ws.Range("B25").Value = "Pernod Ricard"
$ws.Range("C25").Value = "Consumer Defensive"
$ws.Range("D25").Value = "Distillers & Vintners"

$ws.Range("A26").Value = "PUB.PA"
$ws.Range("B26").Value = "Publicis"
$ws.Range("C26").Value = "Communication Services"
$ws.Range("D26").Value = "Advertising"

$ws.Range("A27").Value = "RNO.PA"
$ws.Range("B27").Value = "Renault"
$ws.Range("C27").Value = "Consumer Cyclical"
$ws.Range("D27").Value = "Automobile Manufacturers"

$ws.Range("A28").Value = "SAF.PA"
$ws.Range("B28").Value = "Safran"
$ws.Range("C28").Value = "Industrials"
$ws.Range("D28").Value = "Aerospace & Defense"

$ws.Range("A29").Value = "SGO.PA"
$ws.Range("B29").Value = "Saint-Gobain"
$ws.Range("C29").Value = "Industrials"
$ws.Range("D29").Value = "Building Products"

$ws.Range("A30").Value = "SAN.PA"
$ws.Range("B30").Value = "Sanofi"
$ws.Range("C30").Value = "Healthcare"
$ws.Range("D30").Value = "Pharmaceuticals"

$ws.Range("A31").Value = "SU.PA"
$ws.Range("B31").Value = "Schneider Electric"
$ws.Range("C31").Value = "Industrials"
$ws.Range("D31").Value = "Electrical Components & Equipment"

$ws.Range("A32").Value = "GLE.PA"
$ws.Range("B32").Value = "Société Générale"
$ws.Range("C32").Value = "Financial Services"
$ws.Range("D32").Value = "Diversified Banks"

$ws.Range("A33").Value = "STLAP.PA"
$ws.Range("B33").Value = "Stellantis"
$ws.Range("C33").Value = "Consumer Cyclical"
$ws.Range("D33").Value = "Automobile Manufacturers"

$ws.Range("A34").Value = "STMPA.PA"
$ws.Range("B34").Value = "STMicroelectronics"
$ws.Range("C34").Value = "Technology"
$ws.Range("D34").Value = "Semiconductors"

$ws.Range("A35").Value = "TEP.PA"
$ws.Range("B35").Value = "Teleperformance"
$ws.Range("C35").Value = "Communication Services"
$ws.Range("D35").Value = "Outsourcing"

$ws.Range("A36").Value = "HO.PA"
$ws.Range("B36").Value = "Thales"
$ws.Range("C36").Value = "Industrials"
$ws.Range("D36").Value = "Aerospace & Defense"

$ws.Range("A37").Value = "TTE.PA"
$ws.Range("B37").Value = "TotalEnergies"
$ws.Range("C37").Value = "Energy"
$ws.Range("D37").Value = "Integrated Oil & Gas"

$ws.Range("A38").Value = "URW.PA"
$ws.Range("B38").Value = "Unibail-Rodamco-Westfield"
$ws.Range("C38").Value = "Real Estate"
$ws.Range("D38").Value = "Retail REITs"

$ws.Range("A39").Value = "VIE.PA"
$ws.Range("B39").Value = "Veolia"
$ws.Range("C39").Value = "Industrials"
$ws.Range("D39").Value = "Multi-Utilities"

$ws.Range("A40").Value = "DG.PA"
$ws.Range("B40").Value = "Vinci"
$ws.Range("C40").Value = "Industrials"
$ws.Range("D40").Value = "Construction & Engineering"
